$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Gen slack" (sheet1.xml)
# ---------------------------------------------------------------------------
$wsGenSlack = $wb.Worksheets.Item("Gen slack")
$wsGenSlack.Range("B3").Value = 0
$wsGenSlack.Range("B4").Value = 0
$wsGenSlack.Range("F9").Select()

# ---------------------------------------------------------------------------
# Sheet "Bus" (sheet2.xml)
# ---------------------------------------------------------------------------
$wsBus = $wb.Worksheets.Item("Bus")
$wsBus.Range("B4").Value = 33
$wsBus.Range("B5").Value = 33
$wsBus.Range("B6").Value = 11
$wsBus.Range("B7").Value = 30
$wsBus.Range("B8").Value = 30
$wsBus.Range("B9").Value = 11
$wsBus.Range("B10").Value = 11
$wsBus.Rows("11:12").Delete()
$wsBus.Range("F22").Select()

# ---------------------------------------------------------------------------
# Sheet "Trans" is untouched by this edit.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Sheet "Lines" (sheet4.xml)
# ---------------------------------------------------------------------------
$wsLines = $wb.Worksheets.Item("Lines")
$wsLines.Range("C2").Value = 8
$wsLines.Range("D2").Value = 125
$wsLines.Range("C4").Value = 8
$wsLines.Range("D4").Value = 65
$wsLines.Range("D5").Value = 60
$wsLines.Range("B6").Value = 7
$wsLines.Range("C6").Value = 8
$wsLines.Range("D6").Value = 65
$wsLines.Rows("7:11").Select()
$wsLines.Rows("7:11").Delete()

# ---------------------------------------------------------------------------
# Sheet "Load" (sheet5.xml) - becomes the active/selected tab
# ---------------------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("Load")
$wsLoad.Range("B2").Value = 8
$wsLoad.Activate()
$wsLoad.Range("F9").Select()
